$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values (column D) are stored as text; force text format so
# exact string representation (trailing zeros, sig figs) is preserved.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.62"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.372"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05962"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.395"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.482"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8071"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9092"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1416"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07408"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03313"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03067"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09333"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.856"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001580"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04503"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005932"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006072"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005023"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0009794"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.00007789"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.138"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03884"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006075"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1069"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002700"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007177"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005186"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005792"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.044"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002257"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001997"
